$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data for the newly populated rows 22-25 (missions 20-23)
$data = @(
    @{ row = 22; A = 20; B = 3.26; C = 9.8800000000000008; D = 105;
       E = "launcher and pier";
       F = "clusters from the end of pier and a cross shore line from pier";
       G = "Christine Baker"; H = "EJ Rainville"; I = "Alex de Klerk"; J = "Jim Thomson"; K = "EJ Rainville";
       L = 51; M = 5;
       N = "54,2,3,4,5,6,7,8,9,56,11,12,13,14,57,16,17,18,19,20,21,22,23,24,58,26,27,28,29,59,31,32,33,34,35,36,37,38,40,41,42,43,44,45,46,47,48,49,39,50";
       O = "2,3,4,6,7,8,9,56,11,12,13,14,57,16,17,18,20,21,22,23,24,58,26,27,28,29,59,31,32,33,34,35,36,37,38,40,41,42,43,44,45,46,48,49,39,50";
       P = "1,3,4,5,6";
       Q = "2021-10-10T13:00:00";
       R = "2021-10-10T13:27:00" },
    @{ row = 23; A = 21; B = 3.26; C = 9.8800000000000008; D = 105;
       E = "launcher and pier";
       F = "clusters from the end of pier and a cross shore line from pier";
       G = "Christine Baker"; H = "EJ Rainville"; I = "Alex de Klerk"; J = "Jim Thomson"; K = "EJ Rainville";
       L = 51; M = 5;
       N = "54,2,3,4,6,7,8,9,56,11,12,13,14,57,16,17,18,19,20,21,22,23,24,58,26,27,28,29,59,31,32,33,34,35,36,37,38,40,41,42,43,44,45,46,47,48,49,39,50";
       O = "2,3,4,6,7,8,9,56,11,12,13,14,57,16,17,18,20,21,22,23,24,58,26,27,28,29,59,31,32,33,34,35,36,37,38,40,41,42,43,44,45,46,48,49,39,50";
       P = "1,3,4,5,6";
       Q = "2021-10-10T14:00:00";
       R = "2021-10-10T14:36:00" },
    @{ row = 24; A = 22; B = 3.26; C = 9.8800000000000008; D = 105;
       E = "launcher and pier";
       F = "clusters from the end of pier and a cross shore line from pier";
       G = "Christine Baker"; H = "EJ Rainville"; I = "Alex de Klerk"; J = "Jim Thomson"; K = "EJ Rainville";
       L = 47; M = 5;
       N = "2,3,4,6,7,8,9,56,11,12,13,14,57,16,17,18,20,21,22,23,24,58,26,27,28,29,59,31,32,33,34,35,36,37,38,40,41,42,43,44,45,46,48,49,39,50";
       O = "2,3,4,6,7,8,9,56,11,12,13,14,57,16,17,18,20,21,22,23,24,58,26,27,28,29,59,31,32,33,34,35,36,37,38,40,41,42,43,44,45,46,48,49,39,50";
       P = "1,3,4,5,6";
       Q = "2021-10-10T15:00:00";
       R = "2021-10-10T15:27:00" },
    @{ row = 25; A = 23; B = 3.26; C = 9.8800000000000008; D = 105;
       E = "launcher and pier";
       F = "clusters from the end of pier and a cross shore line from pier";
       G = "Christine Baker"; H = "EJ Rainville"; I = "Alex de Klerk"; J = "Jim Thomson"; K = "EJ Rainville";
       L = 47; M = 5;
       N = "2,3,4,6,7,8,9,56,11,12,13,14,57,16,17,18,20,21,22,23,24,58,26,27,28,29,59,31,32,33,34,35,36,37,38,40,41,42,43,44,45,46,48,49,39,50";
       O = "2,3,4,6,7,8,9,56,11,12,13,14,57,16,17,18,20,21,22,23,24,58,26,27,28,29,59,31,32,33,34,35,36,37,38,40,41,42,43,44,45,46,48,49,39,50";
       P = "1,3,4,5,6";
       Q = "2021-10-10T16:00:00";
       R = "2021-10-10T16:35:00" }
)

foreach ($entry in $data) {
    $r = $entry.row
    $ws.Cells.Item($r, 1).Value = $entry.A
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 6).Value = $entry.F
    $ws.Cells.Item($r, 7).Value = $entry.G
    $ws.Cells.Item($r, 8).Value = $entry.H
    $ws.Cells.Item($r, 9).Value = $entry.I
    $ws.Cells.Item($r, 10).Value = $entry.J
    $ws.Cells.Item($r, 11).Value = $entry.K
    $ws.Cells.Item($r, 12).Value = $entry.L
    $ws.Cells.Item($r, 13).Value = $entry.M
    $ws.Cells.Item($r, 14).Value = $entry.N
    $ws.Cells.Item($r, 15).Value = $entry.O
    $ws.Cells.Item($r, 16).Value = $entry.P
    $ws.Cells.Item($r, 17).Value = $entry.Q
    $ws.Cells.Item($r, 18).Value = $entry.R
    $ws.Rows.Item($r).RowHeight = 85
}

# Update view: selection (and, best-effort, scroll position) to mirror the
# commit's screen state. Selecting the multi-cell range B25:K25 matches the
# saved <selection activeCell="B25" sqref="B25:K25"/>.
$ws.Range("B25:K25").Select()
try { $excel.ActiveWindow.ScrollRow = 17 } catch { }
